# Atualização de bases das ligas, do dia: 14-06-2024 às 20:31
# This script swaps the match data (everything except id/Div/Date columns)
# between two pairs of adjacent rows: (114,115) and (173,174).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($Row1, $Row2) {
    # Column B (2) through AD (30) hold the match data that gets swapped.
    # Columns A (1), C (3), D (4) are left untouched.
    for ($col = 2; $col -le 30; $col++) {
        if ($col -eq 3 -or $col -eq 4) {
            continue
        }
        $cell1 = $ws.Cells.Item($Row1, $col)
        $cell2 = $ws.Cells.Item($Row2, $col)

        $val1 = $cell1.Value2
        $val2 = $cell2.Value2

        $cell1.Value2 = $val2
        $cell2.Value2 = $val1
    }
}

Swap-RowData 114 115
Swap-RowData 173 174
